$p = $ppt.ActivePresentation

# Slide 1: Title "First slide" -> split the leading "First " run into
# two runs: "First" and " " (so the trailing space becomes its own run,
# matching the un-consolidated run layout).
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$space1 = $tr1.Characters(6, 1)
$space1.Text = " "

# Slide 3: Title "Third slide" -> same split for the leading "Third " run.
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$space3 = $tr3.Characters(6, 1)
$space3.Text = " "
